$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 3.781711156805759

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 14.36450238910742

$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 3099.503889238888
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 3112.200597044728
